$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

Set-TextValue 'D2' '60.146.47'
Set-TextValue 'E2' '  -1.31%  '
Set-TextValue 'D3' '2.626.60'
Set-TextValue 'E3' '  +0.95%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '518.79'
Set-TextValue 'E5' '  -1.02%  '
Set-TextValue 'D6' '147.84'
Set-TextValue 'E6' '  -4.50%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.12%  '
Set-TextValue 'E8' '  -3.75%  '
Set-TextValue 'D9' '2.634.63'
Set-TextValue 'E9' '  +0.85%  '
Set-TextValue 'D10' '6.32'
Set-TextValue 'E10' '  -5.63%  '
Set-TextValue 'E11' '  -0.88%  '
Set-TextValue 'E12' '  -2.32%  '
Set-TextValue 'E13' '  -0.69%  '
Set-TextValue 'D14' '3.086.65'
Set-TextValue 'E14' '  +0.90%  '
Set-TextValue 'D15' '60.152.23'
Set-TextValue 'E15' '  -1.34%  '
Set-TextValue 'D16' '21.17'
Set-TextValue 'E16' '  -2.60%  '
Set-TextValue 'E17' '  -2.05%  '
Set-TextValue 'D18' '2.629.57'
Set-TextValue 'E18' '  +0.95%  '
Set-TextValue 'E19' '  -2.47%  '
Set-TextValue 'D20' '340.51'
Set-TextValue 'E20' '  -3.69%  '
Set-TextValue 'D21' '10.40'
Set-TextValue 'E21' '  -1.74%  '
Set-TextValue 'E22' '  -1.70%  '
Set-TextValue 'E23' '  -0.44%  '
Set-TextValue 'D24' '61.16'
Set-TextValue 'E24' '  +0.00%  '
Set-TextValue 'E25' '  -2.30%  '
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  +0.13%  '
Set-TextValue 'D27' '0.159'
Set-TextValue 'E27' '  -4.11%  '
Set-TextValue 'D28' '0.0₃0808'
Set-TextValue 'E28' '  -4.85%  '
Set-TextValue 'D29' '7.03'
Set-TextValue 'E29' '  -4.70%  '
Set-TextValue 'E30' '  +0.00%  '
Set-TextValue 'E31' '  -1.22%  '
Set-TextValue 'E32' '  -5.73%  '
Set-TextValue 'D33' '18.92'
Set-TextValue 'E33' '  -2.43%  '
Set-TextValue 'D34' '150.01'
Set-TextValue 'E34' '  +0.46%  '
Set-TextValue 'E35' '  -7.26%  '
Set-TextValue 'D36' '0.919'
Set-TextValue 'E36' '  -3.22%  '
Set-TextValue 'E37' '  -5.97%  '
Set-TextValue 'E38' '  +0.91%  '
Set-TextValue 'E39' '  +0.63%  '
Set-TextValue 'E40' '  -4.98%  '
Set-TextValue 'E41' '  -4.40%  '
Set-TextValue 'D42' '290.41'
Set-TextValue 'E42' '  +0.95%  '
Set-TextValue 'D43' '0.629'
Set-TextValue 'E43' '  +0.43%  '
Set-TextValue 'D44' '0.0998'
Set-TextValue 'E44' '  -1.42%  '
Set-TextValue 'E45' '  +0.07%  '
Set-TextValue 'E46' '  -2.60%  '
Set-TextValue 'D47' '19.44'
Set-TextValue 'E47' '  -1.00%  '
Set-TextValue 'B48' 'VeChain'
Set-TextValue 'C48' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D48' '0.0232'
Set-TextValue 'E48' '  -2.29%  '
Set-TextValue 'B49' 'WhiteBITCoin'
Set-TextValue 'C49' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D49' '10.39'
Set-TextValue 'E49' '  +0.77%  '
Set-TextValue 'E50' '  -6.93%  '
Set-TextValue 'D51' '1.956.99'
Set-TextValue 'E51' '  -0.17%  '
